$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B5").Value = -1
$ws.Range("C2:C5").Value = 2

$ws.Range("C4").Select()
